# "added work to q7" - swap the BLS series from gasoline prices to
# electricity prices (series APU000072610) and update the 2022 monthly
# values (row 11) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Series Id
$ws.Range("B4").Value = "APU000072610"

# Series Title
$ws.Range("B5").Value = "Electricity per KWH in U.S. city average, average price, not seasonally adjusted"

# Item
$ws.Range("B7").Value = "Electricity per KWH"

# 2022 monthly values (Jan..Dec) for the electricity series
$ws.Range("B11").Value = 0.147
$ws.Range("C11").Value = 0.148
$ws.Range("D11").Value = 0.15
$ws.Range("E11").Value = 0.151
$ws.Range("F11").Value = 0.154
$ws.Range("G11").Value = 0.16
$ws.Range("H11").Value = 0.164
$ws.Range("I11").Value = 0.167
$ws.Range("J11").Value = 0.167
$ws.Range("K11").Value = 0.166
$ws.Range("L11").Value = 0.163
$ws.Range("M11").Value = 0.165
